# Edit: reemplazando caldas por manizales
# Renames the "Caldas" row/column to "Manizales" and updates its toll-cost
# values, while keeping the (now independent) "Medellin" column/row frozen
# at the values it used to mirror from "Caldas".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the city label "Caldas" -> "Manizales" (row 8 label, column H header)
$ws.Range("A8").Value = "Manizales"
$ws.Range("H1").Value = "Manizales"

# 2. New "Manizales" row (row 8) toll-cost values
$ws.Range("B8").Value = 39700
$ws.Range("C8").Value = 117500
$ws.Range("D8").Value = 45400
$ws.Range("E8").Value = 51300
$ws.Range("F8").Value = 136400
$ws.Range("G8").Value = 66400
$ws.Range("I8").Value = 29400
$ws.Range("J8").Value = 106100
$ws.Range("K8").Value = 75800
$ws.Range("L8").Value = 24400
$ws.Range("M8").Value = 104700
$ws.Range("P8").Value = 78700

# 3. New "Manizales" column (H) toll-cost values for the remaining rows
$ws.Range("H2").Value = 39700
$ws.Range("H3").Value = 117500
$ws.Range("H4").Value = 45400
$ws.Range("H5").Value = 51300
$ws.Range("H6").Value = 136400
$ws.Range("H7").Value = 66400
$ws.Range("H9").Value = 29400
$ws.Range("H10").Value = 106100
$ws.Range("H11").Value = 75800
$ws.Range("H12").Value = 24400
$ws.Range("H13").Value = 104700
$ws.Range("H16").Value = 78700

# 4. Freeze the "Medellin" column (I) values that used to mirror column H via
#    formula - they keep their old numeric values as literals now.
$ws.Range("I4").Value = 74800
$ws.Range("I5").Value = 70100
$ws.Range("I6").Value = 107000
$ws.Range("I7").Value = 85200
$ws.Range("I9").Value = 0
$ws.Range("I10").Value = 90700
$ws.Range("I11").Value = 93600
$ws.Range("I12").Value = 42200
$ws.Range("I13").Value = 122500
$ws.Range("I14").Value = 112400
$ws.Range("I15").Value = 93600
$ws.Range("I16").Value = 97500

# 5. Freeze the "Medellin" row (row 9) values that used to mirror row 8 via
#    formula - they keep their old numeric values as literals now.
$ws.Range("D9").Value = 74800
$ws.Range("E9").Value = 70100
$ws.Range("F9").Value = 107000
$ws.Range("G9").Value = 85200
$ws.Range("J9").Value = 90700
$ws.Range("K9").Value = 93600
$ws.Range("L9").Value = 42200
$ws.Range("M9").Value = 122500
$ws.Range("N9").Value = 112400
$ws.Range("O9").Value = 93600
$ws.Range("P9").Value = 97500

# 6. Update the saved selection to match the commit (P9 instead of P14)
$ws.Range("P9").Select()
